$wb = $excel.ActiveWorkbook
$nl = [char]10

# List of (old, new) text replacements to apply to every worksheet's shared
# strings. These remove the footnote markers like " [1]" / " [5]" etc. and
# collapse embedded line-breaks into single spaces (joining multi-line
# labels into one line), matching the author's cleanup pass.
$replacements = @(
    @("DTaP [1]", "DTaP "),
    @("DTaP-IPV [2]", "DTaP-IPV "),
    @("DTaP-Hep B-IPV [4]", "DTaP-Hep B-IPV "),
    @("DTaP-IP-HI [4]", "DTaP-IP-HI "),
    @("e-IPV [5]", "e-IPV "),
    @("Hepatitis A Pediatric [5]", "Hepatitis A Pediatric "),
    @("Hepatitis A-Hepatitis B 18 only [3]", "Hepatitis A-Hepatitis B 18 only "),
    @("Hepatitis B [5]${nl}Pediatric/Adolescent", "Hepatitis B  Pediatric/Adolescent"),
    @("Recombivax${nl}HB", "Recombivax HB"),
    @("Hib [5]", "Hib "),
    @("HPV - Human Papillomavirus 9-valent [5]", "HPV - Human Papillomavirus 9-valent "),
    @("MENB - Meningococcal Group B [5]", "MENB - Meningococcal Group B "),
    @("Meningococcal Conjugate (Groups A, C, Y and W-135) [5]", "Meningococcal Conjugate (Groups A, C, Y and W-135) "),
    @("Measles, Mumps and Rubella (MMR) [1]", "Measles, Mumps and Rubella (MMR) "),
    @("MMR/Varicella [2]", "MMR/Varicella "),
    @("Pneumococcal${nl}13-valent [5] (Pediatric)", "Pneumococcal 13-valent  (Pediatric)"),
    @("Rotavirus, Live, Oral, Pentavalent [5]", "Rotavirus, Live, Oral, Pentavalent "),
    @("Rotavirus, Live, Oral, Oral [5]", "Rotavirus, Live, Oral, Oral "),
    @("Tetanus and Diphtheria Toxoids [3]", "Tetanus and Diphtheria Toxoids "),
    @("Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]", "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis "),
    @("Varicella [5]", "Varicella "),
    @("Hepatitis A Adult [5]", "Hepatitis A Adult "),
    @("Hepatitis A-Hepatitis B Adult [3]", "Hepatitis A-Hepatitis B Adult "),
    @("Hepatitis B Adult [5]", "Hepatitis B Adult "),
    @("HPV-Human Papillomavirus 9 Valent [5]", "HPV-Human Papillomavirus 9 Valent "),
    @("Measles, Mumps,  Rubella [1]", "Measles, Mumps,  Rubella "),
    @("Pneumococcal${nl}13-valent [5]", "Pneumococcal 13-valent "),
    @("Influenza [5]${nl}(Age 6 months and older)", "Influenza  (Age 6 months and older)"),
    # NOTE: the 3-line "Fluzone/Quadrivalent/Pediatric dose" string must be
    # replaced before the shorter 2-line "Fluzone/Quadrivalent" string,
    # since the latter is a prefix of the former.
    @("Fluzone${nl}Quadrivalent${nl}Pediatric dose", "Fluzone Quadrivalent Pediatric dose"),
    @("Fluzone${nl}Quadrivalent", "Fluzone Quadrivalent"),
    @("Influenza [5]${nl}(Age 6-35 months)", "Influenza  (Age 6-35 months)"),
    @("Influenza [5]${nl}(Age 36 months and older)", "Influenza  (Age 36 months and older)"),
    @("Fluarix${nl}Quadrivalent", "Fluarix Quadrivalent"),
    @("FluLaval${nl}Quadrivalent", "FluLaval Quadrivalent"),
    @("Influenza [5]${nl}(Age 4 years and older)", "Influenza  (Age 4 years and older)"),
    @("Influenza [5]${nl}(Age 5 years and older)", "Influenza  (Age 5 years and older)"),
    @("Influenza [5]${nl}Live, Intranasal (Age 2-49 years)", "Influenza  Live, Intranasal (Age 2-49 years)"),
    @("FluMist${nl}Quadrivalent", "FluMist Quadrivalent"),
    @("Afluria${nl}Quadrivalent", "Afluria Quadrivalent")
)

foreach ($ws in $wb.Worksheets) {
    foreach ($pair in $replacements) {
        [void]$ws.Cells.Replace($pair[0], $pair[1])
    }
}
